$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet "Parametros Pol" at the end (after the last existing sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Parametros Pol"

# Copy the header cell format (bold, centered, bordered) from an existing
# formatted cell so we reuse the existing style instead of creating new ones
$ws1.Range("A1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)
$ws.Range("A1").ClearContents()

# Row 1 header / spacer labels
$ws.Range("B1").Value = " "
$ws.Range("C1").Value = "  "
$ws.Range("D1").Value = "   "

# Column A
$ws.Range("A2").Value = "ΣOBS"
$ws.Range("A3").Value = "ΣTeo"
$ws.Range("A4").Value = "e ang"
$ws.Range("A5").Value = "e perm"
$ws.Range("A6").Value = "corr ang"

# Column B
$ws.Range("B2").Value = "900° 0'20.0"
$ws.Range("B3").Value = "900° 0'0.0"
$ws.Range("B4").Value = "0° 0'20.0"
$ws.Range("B5").Value = "0° 0'10.0"
$ws.Range("B6").Value = "-0° 0'5.0"

# Column C
$ws.Range("C2").Value = "ΣDIST"
$ws.Range("C3").Value = "ΔPNS"
$ws.Range("C4").Value = "ΔPEW"
$ws.Range("C5").Value = "e dist"
$ws.Range("C6").Value = "P"

# Column D (numeric values)
$ws.Range("D2").Value = 222.807
$ws.Range("D3").Value = -0.01407894390836262
$ws.Range("D4").Value = 0.01240938702116523
$ws.Range("D5").Value = 0.01876724667648094
$ws.Range("D6").Value = 11872.11975421099

$ws1.Activate()
